# Apply the BOM update described by the commit:
# "Added a line to motor controller bom, added both libraries, a schematic, and a pcb file"
#
# Workbook-level effect visible in the diff: a new BOM line is appended as row 23
# (columns A and C populated), while the pre-existing "Total" row values in F23/G23
# are preserved. This also grows sharedStrings.xml by the two new string values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM entry text added on row 23 (Total row F23/G23 stays where it is).
$ws.Range("A23").Value = "NEW ITEM"
$ws.Range("C23").Value = "FOR TESTING CONFLICT RESOLUTION SKILLS OF GITHUB"

# Reflect the updated view/selection state from the saved workbook: the window had
# scrolled down (top-left visible cell around A10) and the active selection moved
# from H12 to C25.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C25").Select()
